$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 26-34 (shift cells up) to shrink dimension to F25
$ws.Rows("26:34").Delete() | Out-Null

# D2, E2, F2, F3 special updates
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "NSE:LICHSGFIN"
$ws.Range("F2").Value = "NSE:DELHIVERY"
$ws.Range("F3").Value = "NSE:NTPC"

# Column B updates (rows 2-25)
$ws.Range("B2").Value = "NSE:21STCENMGM"
$ws.Range("B3").Value = "NSE:ACE"
$ws.Range("B4").Value = "NSE:APTUS"
$ws.Range("B5").Value = "NSE:ARIES"
$ws.Range("B6").Value = "NSE:BECTORFOOD"
$ws.Range("B7").Value = "NSE:CUMMINSIND"
$ws.Range("B8").Value = "NSE:CYIENTDLM"
$ws.Range("B9").Value = "NSE:DBREALTY"
$ws.Range("B10").Value = "NSE:DCM"
$ws.Range("B11").Value = "NSE:DELHIVERY"
$ws.Range("B12").Value = "NSE:DPSCLTD"
$ws.Range("B13").Value = "NSE:GEPIL"
$ws.Range("B14").Value = "NSE:GRMOVER"
$ws.Range("B15").Value = "NSE:HATSUN"
$ws.Range("B16").Value = "NSE:INDTERRAIN"
$ws.Range("B17").Value = "NSE:JUBLINGREA"
$ws.Range("B18").Value = "NSE:KAJARIACER"
$ws.Range("B19").Value = "NSE:KITEX"
$ws.Range("B20").Value = "NSE:KOTAKBANK"
$ws.Range("B21").Value = "NSE:LUXIND"
$ws.Range("B22").Value = "NSE:MANGALAM"
$ws.Range("B23").Value = "NSE:NITINSPIN"
$ws.Range("B24").Value = "NSE:NTPC"
$ws.Range("B25").Value = "NSE:ROTO"

# Column C updates (rows 2-21), clear C22:C25
$ws.Range("C2").Value = "NSE:ABAN"
$ws.Range("C3").Value = "NSE:AJOONI"
$ws.Range("C4").Value = "NSE:ASTERDM"
$ws.Range("C5").Value = "NSE:BLUESTARCO"
$ws.Range("C6").Value = "NSE:BPL"
$ws.Range("C7").Value = "NSE:CLEAN"
$ws.Range("C8").Value = "NSE:DANGEE"
$ws.Range("C9").Value = "NSE:DEVIT"
$ws.Range("C10").Value = "NSE:DREAMFOLKS"
$ws.Range("C11").Value = "NSE:DUCON"
$ws.Range("C12").Value = "NSE:E2E"
$ws.Range("C13").Value = "NSE:GULFPETRO"
$ws.Range("C14").Value = "NSE:INDOCO"
$ws.Range("C15").Value = "NSE:LTFOODS"
$ws.Range("C16").Value = "NSE:MICEL"
$ws.Range("C17").Value = "NSE:MMTC"
$ws.Range("C18").Value = "NSE:MUKANDLTD"
$ws.Range("C19").Value = "NSE:MUNJALAU"
$ws.Range("C20").Value = "NSE:NAVKARCORP"
$ws.Range("C21").Value = "NSE:NETWORK18"
$ws.Range("C22:C25").Value = ""
